# Menu filtre de la boutique
# Adds 4 new journal entries (rows 46-49) to the "Journal de travail" sheet,
# two of which introduce brand-new description strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 42 is a "plain" existing data row (no special row height, no one-off
# cell style) - use it as the formatting template for the new rows so the
# new cells inherit the same number formats / borders / alignment.
$templateRow = $ws.Range("A42:F42")

foreach ($r in 46..49) {
    $templateRow.Copy() | Out-Null
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Row 46: 2023-05-22 (serial 45068), wraps up the display-bug fixes
$ws.Cells.Item(46, 1).Value = 45068
$ws.Cells.Item(46, 2).Value = 4
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(46, 4).Value = "Implémentation"
$ws.Cells.Item(46, 5).Value = "Corriger les petites erreurs d'affichage"

# Row 47: 2023-05-23 (serial 45069), new shop filter-menu task
$ws.Cells.Item(47, 1).Value = 45069
$ws.Cells.Item(47, 2).Value = 4
$ws.Cells.Item(47, 3).Value = 1.5
$ws.Cells.Item(47, 4).Value = "Implémentation"
$ws.Cells.Item(47, 5).Value = "Créer un menu de filtre sur la page boutique"

# Row 48: 2023-05-23, new cart-checkout task (wraps to two lines like row 41/45)
$ws.Cells.Item(48, 1).Value = 45069
$ws.Cells.Item(48, 2).Value = 4
$ws.Cells.Item(48, 3).Value = 1.5
$ws.Cells.Item(48, 4).Value = "Implémentation"
$ws.Cells.Item(48, 5).Value = "Valider les commandes du panier et envoyer les données dans la base de donées "
$ws.Rows.Item(48).RowHeight = 30

# Row 49: 2023-05-23, documentation task
$ws.Cells.Item(49, 1).Value = 45069
$ws.Cells.Item(49, 2).Value = 4
$ws.Cells.Item(49, 3).Value = 2.25
$ws.Cells.Item(49, 4).Value = "Documentation"
$ws.Cells.Item(49, 5).Value = "Avancer sur la documentation"

# Update the sheet's view state to mirror the saved file: scrolled down a
# few more rows, selection now sitting on the new last row's F column.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("F50").Select()
